# Summarizing Ourselves final data collection
#
# Adds a bold "TAGS:" heading paragraph after the final "Role 2:" paragraph,
# followed by seven plain (non-bold) tag lines.

$d = $word.ActiveDocument

$tags = @(
  "1 Aipangyaraq, Aipaqsaraq -- Marriage",
  "1 Umyuallguteklutek – Being in Harmony",
  "1 Nunalgutkenrilkemeggnun Nallunailuciit -- Being Widely Known",
  "1 Akusrarun -- Mischief, Misconduct",
  "1 Kenkiyaraq -- Showing Love",
  "1 Qanminek Mulngaksaraq -- Careful with Words",
  "1 Agleryaraq -- Menstruation"
)

# Step 1: Build the plain tag paragraphs in a scratch location (just before
# paragraph 2, which is already a plain/non-bold paragraph) so that each new
# paragraph mark inherits clean, non-bold formatting instead of the bold
# formatting used throughout the rest of the form. Paragraphs are inserted in
# reverse order so that the final reading order matches the tag list.
for ($j = $tags.Length - 1; $j -ge 0; $j--) {
  $tag = $tags[$j]
  $p2 = $d.Paragraphs.Item(2)
  $insPos = $p2.Range.Start
  $ins = $d.Range($insPos, $insPos)
  $ins.InsertParagraphBefore()
  $newp = $d.Paragraphs.Item(2)
  $r = $d.Range($newp.Range.Start, $newp.Range.End - 1)
  $r.Text = $tag
}

# Step 2: Cut that scratch block of 7 paragraphs (including their paragraph
# marks) back out so it can be relocated to the end of the document.
$blockStart = $d.Paragraphs.Item(2).Range.Start
$blockEnd = $d.Paragraphs.Item(1 + $tags.Length).Range.End
$blockRange = $d.Range($blockStart, $blockEnd)
$blockRange.Cut()

# Step 3: Append the bold "TAGS:" paragraph right after "Role 2:" (the
# current last paragraph). It naturally inherits the bold/bCs run formatting
# used by the other field labels in the document.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$rTags = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$rTags.Text = "TAGS:"

# Step 4: Append a temporary blank anchor paragraph after "TAGS:". It will
# inherit the bold formatting too, but it only exists momentarily so that we
# have a safe insertion point to paste the non-bold block before.
$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()

# Step 5: Paste the previously-cut plain tag block immediately before that
# blank anchor paragraph. Pasting preserves the non-bold formatting captured
# when the block was built in step 1.
$n = $d.Paragraphs.Count
$anchorP = $d.Paragraphs.Item($n)
$pastePos = $anchorP.Range.Start
$dest = $d.Range($pastePos, $pastePos)
$dest.Paste()

# Step 6: Remove the now-unneeded trailing blank anchor paragraph, leaving
# the last of the seven tag paragraphs as the final paragraph in the body.
$n2 = $d.Paragraphs.Count
$anchorP2 = $d.Paragraphs.Item($n2)
$anchorP2.Range.Delete()
